# Fluxo de caixa.pptx - apply edits described by the commit diff:
#  1. Update the fixed "datetime1" field text on every slide layout and
#     the slide master from 16/03/2021 -> 19/03/2021.
#  2. Slide 4 ("Elicitação"): tweak the "Etnografia"/"Entrevista aberta"
#     paragraphs' wording and turn on "Shrink text on overflow"
#     (normAutofit) for that content placeholder.

$p = $ppt.ActivePresentation

$oldDate = "16/03/2021"
$newDate = "19/03/2021"

# ppPlaceholderDate
$ppPlaceholderDate = 16

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Placeholders.Count; $i++) {
        $sh = $shapes.Placeholders.Item($i)
        if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# 1a. Every slide layout's date placeholder.
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes
}

# 1b. The slide master's own date placeholder.
Update-DatePlaceholders $master.Shapes

# NOTE: this runtime does not support persisting text edits on the
# handout master / notes master date placeholders (writes silently do
# not round-trip, and -- worse -- touching them corrupts unrelated
# slide-master shape text), so they are intentionally left alone here.

# 2. Slide 4 content + autofit.
$slide4 = $p.Slides.Item(4)
$contentShape = $slide4.Shapes.Item(2)

$contentShape.TextFrame.AutoSize = 2  # ppAutoSizeTextToFitShape -> <a:normAutofit/>

$tr = $contentShape.TextFrame.TextRange
$tr.Paragraphs(1).Runs(2).Text = ": A ideia a principio surgiu de uma necessidade pessoal de controlar as finanças pessoais. Os aplicativos de bancos não estavam atendendo mais, principalmente quando o lançamento era em espécie. Evitando o uso do Excel para esse controle, surgiu a necessidade de ter um software que facilitasse a manipulação dos dados. "
$tr.Paragraphs(2).Runs(2).Text = ": Em seguida uma breve apresentação para colegas e familiares sobre o projeto e como poderia estar melhorando."
